$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.975.71'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '2.968.59'
$ws.Range('E3').Value = '  +2.73%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''353.66'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').Value = '''112.26'
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('D7').Value = '''0.566'
$ws.Range('E7').Value = '  +0.97%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('D10').Value = '''39.75'
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('D11').Value = '''0.0898'
$ws.Range('E11').Value = '  +5.13%  '
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').Value = '''19.97'
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('D14').Value = '''7.94'
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').Value = '3.440.40'
$ws.Range('E15').Value = '  +3.11%  '
$ws.Range('D16').Value = '2.981.40'
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('D17').Value = '''0.998'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '52.080.00'
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('D19').Value = '''7.72'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('D20').Value = '''14.57'
$ws.Range('E20').Value = '  +6.60%  '
$ws.Range('E21').Value = '  -2.40%  '
$ws.Range('D22').Value = '0.0₃0992'
$ws.Range('E22').Value = '  +1.28%  '
$ws.Range('D23').Value = '''71.37'
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').Value = '''270.89'
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('D26').Value = '''0.181'
$ws.Range('E26').Value = '  +9.52%  '
$ws.Range('D27').Value = '''27.62'
$ws.Range('E27').Value = '  +3.76%  '
$ws.Range('D28').Value = '''7.65'
$ws.Range('E28').Value = '  +20.71%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '''0.109'
$ws.Range('E30').Value = '  +21.57%  '
$ws.Range('D31').Value = '''10.75'
$ws.Range('E31').Value = '  +1.37%  '
$ws.Range('D32').Value = '''37.78'
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('D33').Value = '''6.22'
$ws.Range('E33').Value = '  +10.17%  '
$ws.Range('D34').Value = '''52.99'
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').Value = '''2.07'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '''0.0450'
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').Value = '''3.44'
$ws.Range('E38').Value = '  +3.08%  '
$ws.Range('D39').Value = '''18.99'
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('D41').Value = '''2.69'
$ws.Range('E41').Value = '  +3.28%  '
$ws.Range('D42').Value = '''23.99'
$ws.Range('E42').Value = '  +5.30%  '
$ws.Range('E43').Value = '  +1.61%  '
$ws.Range('D44').Value = '''2.18'
$ws.Range('E44').Value = '  -2.18%  '
$ws.Range('D45').Value = '''3.56'
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('D47').Value = '2.184.24'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = '''113.83'
$ws.Range('E48').Value = '  -7.33%  '
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').Value = '''0.0341'
$ws.Range('E50').Value = '  +6.37%  '
$ws.Range('D51').Value = '''0.939'
$ws.Range('E51').Value = '  -2.77%  '
